{"js": "const replacements = [\n  [\"170\u00f72=85, 0\", \"505\u00f78=63, 1\"],\n  [\"150\u00f78=18, 6\", \"461\u00f73=153, 2\"],\n  [\"980\u00f73=326, 2\", \"124\u00f74=31, 0\"],\n  [\"677\u00f79=75, 2\", \"274\u00f74=68, 2\"],\n  [\"244\u00f77=34, 6\", \"512\u00f78=64, 0\"],\n  [\"464\u00f75=92, 4\", \"377\u00f76=62, 5\"],\n  [\"847\u00f76=141, 1\", \"565\u00f75=113, 0\"],\n  [\"601\u00f73=200, 1\", \"480\u00f78=60, 0\"],\n  [\"580\u00f74=145, 0\", \"933\u00f74=233, 1\"],\n  [\"364\u00f77=52, 0\", \"730\u00f78=91, 2\"],\n  [\"991\u00f75=198, 1\", \"737\u00f75=147, 2\"],\n  [\"132\u00f75=26, 2\", \"279\u00f73=93, 0\"],\n  [\"226\u00f74=56, 2\", \"302\u00f74=75, 2\"],\n  [\"505\u00f74=126, 1\", \"386\u00f75=77, 1\"],\n  [\"673\u00f73=224, 1\", \"435\u00f72=217, 1\"],\n  [\"657\u00f75=131, 2\", \"670\u00f75=134, 0\"],\n  [\"101\u00f75=20, 1\", \"812\u00f75=162, 2\"],\n  [\"140\u00f76=23, 2\", \"156\u00f75=31, 1\"],\n  [\"707\u00f74=176, 3\", \"392\u00f78=49, 0\"],\n  [\"704\u00f75=140, 4\", \"389\u00f73=129, 2\"],\n  [\"800\u00f75=160, 0\", \"571\u00f78=71, 3\"],\n  [\"979\u00f75=195, 4\", \"832\u00f79=92, 4\"],\n  [\"259\u00f77=37, 0\", \"651\u00f79=72, 3\"],\n  [\"535\u00f79=59, 4\", \"234\u00f74=58, 2\"],\n  [\"338\u00f79=37, 5\", \"543\u00f74=135, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"170\u00f72=85, 0\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"505\u00f78=63, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"150\u00f78=18, 6\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"461\u00f73=153, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"980\u00f73=326, 2\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"124\u00f74=31, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"677\u00f79=75, 2\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"274\u00f74=68, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"244\u00f77=34, 6\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"512\u00f78=64, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"464\u00f75=92, 4\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"377\u00f76=62, 5\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"847\u00f76=141, 1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"565\u00f75=113, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"601\u00f73=200, 1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"480\u00f78=60, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"580\u00f74=145, 0\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"933\u00f74=233, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"364\u00f77=52, 0\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"730\u00f78=91, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"991\u00f75=198, 1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"737\u00f75=147, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"132\u00f75=26, 2\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"279\u00f73=93, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"226\u00f74=56, 2\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"302\u00f74=75, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"505\u00f74=126, 1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"386\u00f75=77, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"673\u00f73=224, 1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"435\u00f72=217, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"657\u00f75=131, 2\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"670\u00f75=134, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"101\u00f75=20, 1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"812\u00f75=162, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"140\u00f76=23, 2\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"156\u00f75=31, 1\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"707\u00f74=176, 3\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"392\u00f78=49, 0\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"704\u00f75=140, 4\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"389\u00f73=129, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"800\u00f75=160, 0\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"571\u00f78=71, 3\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"979\u00f75=195, 4\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"832\u00f79=92, 4\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"259\u00f77=37, 0\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"651\u00f79=72, 3\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"535\u00f79=59, 4\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"234\u00f74=58, 2\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"338\u00f79=37, 5\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"543\u00f74=135, 3\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n"}
